$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to text
# so Excel doesn't auto-convert them from their literal string representation.
$textCells = @("D5", "D6", "D8", "D9", "D11", "D15", "D20", "D22", "D23", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D42", "D46", "D47", "D48", "D49", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.113.99'
$ws.Range("E2").Value = '  -0.74%  '
$ws.Range("D3").Value = '1.653.29'
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("E4").Value = '  -0.48%  '
$ws.Range("D5").Value = '218.88'
$ws.Range("E5").Value = '  -0.73%  '
$ws.Range("D6").Value = '0.5258'
$ws.Range("E6").Value = '  -0.90%  '
$ws.Range("E7").Value = '  -0.43%  '
$ws.Range("D8").Value = '0.2670'
$ws.Range("E8").Value = '  +1.00%  '
$ws.Range("D9").Value = '0.06370'
$ws.Range("E9").Value = '  +0.23%  '
$ws.Range("E10").Value = '  -1.50%  '
$ws.Range("D11").Value = '0.07698'
$ws.Range("E11").Value = '  -1.70%  '
$ws.Range("E12").Value = '  +1.56%  '
$ws.Range("D13").Value = '1.686.63'
$ws.Range("E13").Value = '  +1.03%  '
$ws.Range("D14").Value = '1.880.97'
$ws.Range("E14").Value = '  -0.77%  '
$ws.Range("D15").Value = '0.5610'
$ws.Range("E15").Value = '  +0.16%  '
$ws.Range("D16").Value = '0.0₅8243'
$ws.Range("E16").Value = '  +1.42%  '
$ws.Range("E17").Value = '  -0.45%  '
$ws.Range("D18").Value = '26.115.99'
$ws.Range("E18").Value = '  -0.75%  '
$ws.Range("D20").Value = '4.699'
$ws.Range("E20").Value = '  -0.25%  '
$ws.Range("D22").Value = '191.32'
$ws.Range("E22").Value = '  -3.58%  '
$ws.Range("D23").Value = '5.983'
$ws.Range("E23").Value = '  -1.14%  '
$ws.Range("E24").Value = '  -0.50%  '
$ws.Range("D25").Value = '145.93'
$ws.Range("E25").Value = '  -0.55%  '
$ws.Range("D26").Value = '0.1202'
$ws.Range("E26").Value = '  -0.85%  '
$ws.Range("D27").Value = '7.262'
$ws.Range("E28").Value = '  -1.20%  '
$ws.Range("D29").Value = '1.496'
$ws.Range("E29").Value = '  -1.30%  '
$ws.Range("D30").Value = '0.05655'
$ws.Range("E30").Value = '  -3.90%  '
$ws.Range("D31").Value = '1.272'
$ws.Range("E31").Value = '  -1.18%  '
$ws.Range("D32").Value = '3.502'
$ws.Range("E32").Value = '  -0.92%  '
$ws.Range("D33").Value = '3.386'
$ws.Range("E33").Value = '  +2.16%  '
$ws.Range("D34").Value = '1.581'
$ws.Range("E34").Value = '  -1.33%  '
$ws.Range("D35").Value = '2.801'
$ws.Range("E35").Value = '  -0.96%  '
$ws.Range("D36").Value = '0.9473'
$ws.Range("E36").Value = '  -1.34%  '
$ws.Range("D37").Value = '2.406'
$ws.Range("E37").Value = '  -1.02%  '
$ws.Range("D38").Value = '0.5789'
$ws.Range("E38").Value = '  -0.22%  '
$ws.Range("E39").Value = '  -1.35%  '
$ws.Range("D40").Value = '5.979'
$ws.Range("E40").Value = '  +0.28%  '
$ws.Range("E41").Value = '  -0.51%  '
$ws.Range("D42").Value = '0.8411'
$ws.Range("E42").Value = '  -1.95%  '
$ws.Range("D43").Value = '1.025.75'
$ws.Range("E43").Value = '  -4.45%  '
$ws.Range("E44").Value = '  -1.01%  '
$ws.Range("D45").Value = '1.792.27'
$ws.Range("E45").Value = '  -0.77%  '
$ws.Range("D46").Value = '58.58'
$ws.Range("E46").Value = '  +0.27%  '
$ws.Range("D47").Value = '1.004'
$ws.Range("E47").Value = '  -0.87%  '
$ws.Range("D48").Value = '0.05336'
$ws.Range("E48").Value = '  +3.63%  '
$ws.Range("D49").Value = '8.055'
$ws.Range("E49").Value = '  -0.36%  '
$ws.Range("E50").Value = '  -1.40%  '
$ws.Range("D51").Value = '0.4342'
